$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jul 2022")

# Add a new employee row (row 95) for the "Test Import" record, matching the
# existing table layout: Employee Number, Employee Name, Date of Joining,
# Job Title, Department, Date of birth.
$ws.Cells.Item(95, 1).Value = 189
$ws.Cells.Item(95, 2).Value = "Test Import"
$ws.Cells.Item(95, 3).Value = (Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(95, 4).Value = "abc"
$ws.Cells.Item(95, 5).Value = "abc"
$ws.Cells.Item(95, 6).Value = (Get-Date -Year 1998 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0)

# Reuse the same date-number-format cell style already used in the table
# (copy formats only, so neither value nor any other attribute is touched).
$ws.Range("C94").Copy()
$ws.Range("C95").PasteSpecial(-4122)
$ws.Range("F94").Copy()
$ws.Range("F95").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection/scroll position to mirror the post-edit view.
$ws.Range("A85").Select()
$ws.Range("C95").Select()
